# Saldo.xlsx update:
#  - remove the CLERIA (account 004855960) row entirely
#  - correct LEVI's (account 005206566) balance from 130884.56 to 103884.56

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CLERIA row sits at row 5 (Conta=004855960, Nome=CLERIA, Saldo=187794.6).
# Deleting the whole row shifts everything below it up by one, so LEVI's
# row (originally row 6) becomes row 5.
$ws.Rows.Item(5).Delete()

# Update LEVI's balance.
$ws.Range("C5").Value = 103884.56
